$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

# New requirement rows to append below the existing data (rows 3-11 already exist).
# Columns: B=RequirementID, C=Short Description, D=Requirement Status,
#          E=Developers Assigned, F=Developers Who Completed,
#          G=Standalone/Dependent, H=Priority, I=Critical/Not Critical
$rows = @(
    @{ Id = 10; Desc = "Create Shopping Cart";                      Status = "Complete"; Dev = "Joseph"; Priority = "Standalone"; PriorityNum = 10; Crit = "Critical" },
    @{ Id = 11; Desc = "Edit Shopping Cart";                        Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 10; Crit = "Critical" },
    @{ Id = 12; Desc = "Delete Shopping Cart";                      Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 10; Crit = "Critical" },
    @{ Id = 13; Desc = "Display Shopping Carts";                    Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 10; Crit = "Critical" },
    @{ Id = 14; Desc = "Sign Up as a New User";                     Status = "Complete"; Dev = "Joseph"; Priority = "Standalone"; PriorityNum = 10; Crit = "Critical" },
    @{ Id = 15; Desc = "Login as an Existing User";                 Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 10; Crit = "Critical" },
    @{ Id = 16; Desc = "Display Recipes&Carts on Profile Page";     Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 7;  Crit = "Not Critical" },
    @{ Id = 17; Desc = "Persistent Data between Logins";            Status = "Complete"; Dev = "Joseph"; Priority = "Dependent";  PriorityNum = 10; Crit = "Critical" }
)

$r = 12
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.Id
    $ws.Cells.Item($r, 3).Value = $row.Desc
    $ws.Cells.Item($r, 4).Value = $row.Status
    $ws.Cells.Item($r, 5).Value = $row.Dev
    $ws.Cells.Item($r, 6).Value = $row.Dev
    $ws.Cells.Item($r, 7).Value = $row.Priority
    $ws.Cells.Item($r, 8).Value = $row.PriorityNum
    $ws.Cells.Item($r, 9).Value = $row.Crit

    $ws.Cells.Item($r, 2).HorizontalAlignment = -4131
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131

    $r = $r + 1
}

$ws.Range("I19").Select()
